# Apply cryptos list update (prices / 1h volume % changes, and two row swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.841.39'
$ws.Range('E2').Value = '  +0.78%  '
$ws.Range('D3').Value = '2.491.77'
$ws.Range('E3').Value = '  +1.65%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'" + '532.68'
$ws.Range('E5').Value = '  +0.84%  '
$ws.Range('D6').Value = "'" + '135.79'
$ws.Range('E6').Value = '  +1.41%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  +1.36%  '
$ws.Range('E9').Value = '  +1.77%  '
$ws.Range('E10').Value = '  -1.70%  '
$ws.Range('D11').Value = "'" + '5.42'
$ws.Range('E11').Value = '  +2.34%  '
$ws.Range('E12').Value = '  +1.23%  '
$ws.Range('D13').Value = '2.934.06'
$ws.Range('E13').Value = '  +1.50%  '
$ws.Range('D14').Value = '58.737.45'
$ws.Range('E14').Value = '  +0.75%  '
$ws.Range('D15').Value = "'" + '22.66'
$ws.Range('E15').Value = '  +0.18%  '
$ws.Range('E16').Value = '  +0.09%  '
$ws.Range('D17').Value = '2.490.72'
$ws.Range('E17').Value = '  +1.21%  '
$ws.Range('E18').Value = '  +2.48%  '
$ws.Range('D19').Value = "'" + '4.23'
$ws.Range('E19').Value = '  +0.80%  '
$ws.Range('D20').Value = "'" + '321.79'
$ws.Range('E20').Value = '  +0.41%  '
$ws.Range('E21').Value = '  +0.26%  '
$ws.Range('D22').Value = "'" + '5.96'
$ws.Range('E22').Value = '  +4.19%  '
$ws.Range('D23').Value = "'" + '65.07'
$ws.Range('E23').Value = '  +4.22%  '
$ws.Range('D24').Value = "'" + '0.418'
$ws.Range('E24').Value = '  +2.56%  '
$ws.Range('E25').Value = '  +0.24%  '
$ws.Range('E26').Value = '  +1.38%  '
$ws.Range('E27').Value = '  +0.67%  '
$ws.Range('D28').Value = '0.0₃0757'
$ws.Range('E28').Value = '  +0.65%  '
$ws.Range('E29').Value = '  +4.26%  '
$ws.Range('D30').Value = "'" + '6.43'
$ws.Range('E30').Value = '  -1.02%  '
$ws.Range('E31').Value = '  -0.64%  '
$ws.Range('E32').Value = '  +4.75%  '
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('D34').Value = "'" + '18.29'
$ws.Range('E34').Value = '  +0.41%  '
$ws.Range('E35').Value = '  -0.40%  '
$ws.Range('D36').Value = "'" + '4.00'
$ws.Range('E36').Value = '  -0.29%  '
$ws.Range('D37').Value = "'" + '1.52'
$ws.Range('E37').Value = '  -0.82%  '
$ws.Range('B38').Value = 'Filecoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D38').Value = "'" + '3.55'
$ws.Range('E38').Value = '  +0.36%  '
$ws.Range('B39').Value = 'SuiNetwork'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D39').Value = "'" + '0.792'
$ws.Range('E39').Value = '  -1.57%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').Value = "'" + '5.20'
$ws.Range('E40').Value = '  +2.28%  '
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D41').Value = "'" + '280.32'
$ws.Range('E41').Value = '  +2.27%  '
$ws.Range('E42').Value = '  -0.05%  '
$ws.Range('E43').Value = '  +3.09%  '
$ws.Range('D44').Value = "'" + '129.34'
$ws.Range('E44').Value = '  +7.10%  '
$ws.Range('E45').Value = '  +0.60%  '
$ws.Range('E46').Value = '  -0.26%  '
$ws.Range('E47').Value = '  -1.02%  '
$ws.Range('E48').Value = '  -0.40%  '
$ws.Range('D49').Value = "'" + '17.15'
$ws.Range('E49').Value = '  +0.40%  '
$ws.Range('D50').Value = '1.750.00'
$ws.Range('E50').Value = '  +0.58%  '
$ws.Range('D51').Value = "'" + '0.979'
$ws.Range('E51').Value = '  +0.20%  '
